$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.675.53"
$ws.Range("E2").Value = "  -0.83%  "
$ws.Range("D3").Value = "1.584.81"
$ws.Range("E3").Value = "  -3.10%  "
$ws.Range("E4").Value = "  +0.25%  "
$ws.Range("D5").Value = "206.35"
$ws.Range("E5").Value = "  -2.51%  "
$ws.Range("E6").Value = "  -2.82%  "
$ws.Range("E7").Value = "  +0.27%  "
$ws.Range("D8").Value = "22.24"
$ws.Range("E8").Value = "  -4.73%  "
$ws.Range("E10").Value = "  -3.10%  "
$ws.Range("D11").Value = "0.0867"
$ws.Range("E11").Value = "  -1.83%  "
$ws.Range("D12").Value = "1.810.92"
$ws.Range("E12").Value = "  -3.02%  "
$ws.Range("D13").Value = "1.598.13"
$ws.Range("E13").Value = "  -2.33%  "
$ws.Range("E14").Value = "  -3.97%  "
$ws.Range("D15").Value = "0.531"
$ws.Range("E15").Value = "  -5.55%  "
$ws.Range("D16").Value = "27.656.07"
$ws.Range("E16").Value = "  -0.94%  "
$ws.Range("D17").Value = "63.25"
$ws.Range("E17").Value = "  -3.12%  "
$ws.Range("D18").Value = "220.11"
$ws.Range("E18").Value = "  -3.74%  "
$ws.Range("E19").Value = "  -3.67%  "
$ws.Range("D20").Value = "7.32"
$ws.Range("E20").Value = "  -5.04%  "
$ws.Range("E22").Value = "  -5.03%  "
$ws.Range("D23").Value = "9.49"
$ws.Range("E23").Value = "  -6.37%  "
$ws.Range("D24").Value = "1.96"
$ws.Range("E24").Value = "  -5.84%  "
$ws.Range("D25").Value = "153.92"
$ws.Range("E25").Value = "  -1.33%  "
$ws.Range("B26").Value = "BinanceUSD"
$ws.Range("C26").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.40%  "
$ws.Range("B27").Value = "Cosmos"
$ws.Range("C27").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D27").Value = "6.75"
$ws.Range("E27").Value = "  -2.79%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.10"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.81%  "
$ws.Range("E29").Value = "  -4.06%  "
$ws.Range("E30").Value = "  -2.51%  "
$ws.Range("E31").Value = "  -3.41%  "
$ws.Range("E32").Value = "  -5.62%  "
$ws.Range("D33").Value = "1.385.10"
$ws.Range("E33").Value = "  -1.06%  "
$ws.Range("D34").Value = "2.94"
$ws.Range("E34").Value = "  -5.20%  "
$ws.Range("E35").Value = "  -5.21%  "
$ws.Range("E36").Value = "  -5.23%  "
$ws.Range("E37").Value = "  -1.02%  "
$ws.Range("E38").Value = "  -3.04%  "
$ws.Range("E39").Value = "  -3.17%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.820"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.63%  "
$ws.Range("D42").Value = "0.978"
$ws.Range("E42").Value = "  -2.77%  "
$ws.Range("E43").Value = "  -3.77%  "
$ws.Range("B44").Value = "MXToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D44").Value = "2.17"
$ws.Range("E44").Value = "  +1.69%  "
$ws.Range("B45").Value = "Aave"
$ws.Range("C45").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D45").Value = "63.58"
$ws.Range("E45").Value = "  -3.71%  "
$ws.Range("E46").Value = "  -4.10%  "
$ws.Range("D47").Value = "1.722.37"
$ws.Range("E47").Value = "  -3.02%  "
$ws.Range("D48").Value = "88.12"
$ws.Range("E48").Value = "  -0.62%  "
$ws.Range("E49").Value = "  -2.62%  "
$ws.Range("D50").Value = "0.0974"
$ws.Range("E50").Value = "  -4.92%  "
$ws.Range("D51").Value = "0.0499"
$ws.Range("E51").Value = "  -0.91%  "
